$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the text for the "Valores_consumo Valor" entries (row 11) to "Valores_consumo"
$ws.Range("A11").Value = "Valores_consumo"
$ws.Range("B11").Value = "Valores_consumo"

# Update the active selection to match the saved view state
$ws.Range("F12").Select()
